# Finflux 4723 scenario workbook update
# - Update "Withdrawals blocked..." error message text (sharedStrings content)
# - Move active tab / selection from "NewSavingInput" to "error" sheet
# - Update selection on "error" sheet to B1

$wb = $excel.ActiveWorkbook

# 1) Update the error message text on the "error" sheet (B1).
#    This changes the underlying shared-string text from
#    "Withdrawals blocked until after `01 March 2015`." to
#    "Withdrawals blocked until `01 March 2015`."
$errorSheet = $wb.Worksheets.Item("error")
$errorSheet.Range("B1").Value = "Withdrawals blocked until ``01 March 2015``."

# 2) Move the active sheet from "NewSavingInput" to "error", and update the
#    selected cell on the error sheet to B1 (was B7).
$errorSheet.Activate()
$errorSheet.Range("B1").Select()
